$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The F column ("매크로_상태" / macro status) rows 2-6 previously held
# "미완료" (incomplete) placeholder values; update them to reflect the
# new "오류" (error) status uncovered during the text-search analysis.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 6).Value = "오류"
}
